$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly row at position 160 (data rows are ordered 159 = current,
# then 160.. older dates). This shifts the existing rows 160-183 down to 161-184,
# pushing the former last row (183) into a brand-new row 184.
$ws.Rows.Item(160).Insert()

$ws.Range("A160").Value = 4
$ws.Range("B160").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C160").Value = "Los Lagos"
$ws.Range("D160").Value = 44617
$ws.Range("E160").Value = 10
$ws.Range("F160").Value = 100112039
$ws.Range("G160").Value = "Ciboulette"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 250
$ws.Range("K160").Value = 3000
$ws.Range("L160").Value = 3000
$ws.Range("M160").Value = 3000
$ws.Range("N160").Value = "`$/docena de atados"
$ws.Range("O160").Value = "Región Metropolitana"
$ws.Range("P160").Value = 1000
$ws.Range("Q160").Value = 3
$ws.Range("R160").Value = "Hortaliza"
